# Auto-generated edit script applying the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric strings need to be forced to
# Text format first, otherwise Excel auto-converts them to Number type,
# which would not match the original inline-string ("text") representation.
$textCells = @('D5', 'D6', 'D11', 'D14', 'D18', 'D20', 'D21', 'D23', 'D24', 'D26', 'D27', 'D32', 'D33', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D44', 'D47', 'D48', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '64.029.70'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').Value = '3.138.50'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '591.91'
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').Value = '147.32'
$ws.Range('E6').Value = '  +1.77%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '3.130.86'
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('E10').Value = '  +11.43%  '
$ws.Range('D11').Value = '5.77'
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('E13').Value = '  +3.86%  '
$ws.Range('D14').Value = '37.41'
$ws.Range('E14').Value = '  +4.41%  '
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = '3.652.85'
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '63.878.42'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '7.18'
$ws.Range('E18').Value = '  -2.66%  '
$ws.Range('D19').Value = '3.132.95'
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').Value = '466.96'
$ws.Range('E20').Value = '  +2.04%  '
$ws.Range('D21').Value = '14.37'
$ws.Range('E21').Value = '  +1.30%  '
$ws.Range('E22').Value = '  -0.61%  '
$ws.Range('D23').Value = '7.57'
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('D24').Value = '13.31'
$ws.Range('E24').Value = '  -4.04%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('D27').Value = '9.00'
$ws.Range('E27').Value = '  +7.53%  '
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('E29').Value = '  -1.76%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('D32').Value = '27.19'
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('D33').Value = '0.110'
$ws.Range('E33').Value = '  -4.25%  '
$ws.Range('D34').Value = '0.0₃0890'
$ws.Range('E34').Value = '  +10.03%  '
$ws.Range('E35').Value = '  +7.71%  '
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('D37').Value = '3.41'
$ws.Range('E37').Value = '  +10.71%  '
$ws.Range('D38').Value = '6.13'
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = '457.61'
$ws.Range('E39').Value = '  +6.29%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '50.96'
$ws.Range('E40').Value = '  +0.38%  '
$ws.Range('D41').Value = '8.73'
$ws.Range('E41').Value = '  -2.05%  '
$ws.Range('D42').Value = '0.0374'
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').Value = '2.902.63'
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('D44').Value = '0.279'
$ws.Range('E44').Value = '  -0.58%  '
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('E46').Value = '  -0.95%  '
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').Value = '35.91'
$ws.Range('E47').Value = '  +1.09%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '126.40'
$ws.Range('E48').Value = '  +2.06%  '
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('D51').Value = '24.82'
$ws.Range('E51').Value = '  -0.30%  '
